# SchoolData.xlsx edit: fix up the timeslot table (start/end times) and
# correct one lesson's subject on the lessonList sheet.

$wb = $excel.ActiveWorkbook

# --- timeslotList: correct a couple of start/end times -------------------
$ws1 = $wb.Worksheets.Item("timeslotList")

# Row 1 (MONDAY) started at 08:30:00 - should be 08:31:00
$ws1.Range("C1").Value = "08:31:00"

# Row 2 (MONDAY 09:30-?) and row 7 (TUESDAY 09:30-?) incorrectly ended at
# 11:30:00 instead of 10:30:00
$ws1.Range("D2").Value = "10:30:00"
$ws1.Range("D7").Value = "10:30:00"

# --- lessonList: row 15 (10th grade) was actually Chemistry, not Physics -
$ws3 = $wb.Worksheets.Item("lessonList")
$ws3.Range("B15").Value = "Chemistry"
